# "updated AMGN and added others"
#
# Applies the reachable, substantive edits on the "Main" sheet:
#   - Row 3 (UNH line): refreshed price, quarter tag and as-of date.
#   - Pane selection on the Main sheet moved to B3.
#
# Note: several cells in row 3 (E3/F3/G3/H3) and the row7/row34 external
# formulas are driven by cached values coming from linked external
# workbooks (xl/externalLinks/*.xml). Those caches are refreshed by Excel
# only when the real linked workbooks (UNH.xlsx / SNCE.xlsx / HUM.xlsx)
# are reachable/opened, which isn't something the exposed Application
# object model lets us poke at directly -- they simply recalc from
# whatever is already cached in the file. So here we only touch the
# inputs that are genuinely under the user's control from the UI: the
# unit price, the "as of" quarter label and the "as of" date -- the
# dependent formulas (E3/F3/G3/H3) recompute automatically from them and
# from the existing external caches.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main")

# Row 3 - UNH: new price, quarter tag (Q222 -> Q224) and updated date.
$ws.Range("D3").Value = 572
$ws.Range("I3").Value = "Q224"
$ws.Range("J3").Value2 = 45581

# Move the active selection in the frozen bottom-right pane to B3.
$ws.Activate()
$ws.Range("B3").Select()
